$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 2 4 '67.044.12'
Set-TextValue 2 5 '  -1.99%  '

Set-TextValue 3 4 '3.608.62'
Set-TextValue 3 5 '  -2.71%  '

Set-TextValue 4 4 '1.00'
Set-TextValue 4 5 '  +0.29%  '

Set-TextValue 5 4 '586.82'
Set-TextValue 5 5 '  -2.72%  '

Set-TextValue 6 4 '185.29'
Set-TextValue 6 5 '  +2.30%  '

Set-TextValue 7 4 '0.608'
Set-TextValue 7 5 '  -3.90%  '

Set-TextValue 8 5 '  -0.26%  '

Set-TextValue 9 4 '0.676'
Set-TextValue 9 5 '  -5.76%  '

Set-TextValue 10 4 '0.146'
Set-TextValue 10 5 '  -10.82%  '

Set-TextValue 11 4 '54.22'
Set-TextValue 11 5 '  -4.24%  '

Set-TextValue 12 4 '0.0000252'
Set-TextValue 12 5 '  -13.82%  '

Set-TextValue 13 4 '9.94'
Set-TextValue 13 5 '  -4.78%  '

Set-TextValue 14 4 '4.195.69'
Set-TextValue 14 5 '  -2.32%  '

Set-TextValue 15 4 '3.621.90'
Set-TextValue 15 5 '  -2.27%  '

Set-TextValue 16 5 '  -0.30%  '

Set-TextValue 17 4 '18.41'
Set-TextValue 17 5 '  -4.92%  '

Set-TextValue 18 4 '66.761.35'
Set-TextValue 18 5 '  -2.17%  '

Set-TextValue 19 2 'Uniswap'
Set-TextValue 19 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 19 4 '12.22'
Set-TextValue 19 5 '  -5.26%  '

Set-TextValue 20 2 'Polygon'
Set-TextValue 20 3 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 20 4 '1.07'
Set-TextValue 20 5 '  -4.79%  '

Set-TextValue 21 4 '394.81'
Set-TextValue 21 5 '  -3.73%  '

Set-TextValue 22 4 '4.33'
Set-TextValue 22 5 '  -7.09%  '

Set-TextValue 23 4 '85.48'
Set-TextValue 23 5 '  -4.32%  '

Set-TextValue 24 4 '2.85'
Set-TextValue 24 5 '  -6.03%  '

Set-TextValue 25 4 '12.23'
Set-TextValue 25 5 '  -5.00%  '

Set-TextValue 26 2 'LEO'
Set-TextValue 26 3 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 26 4 '6.05'
Set-TextValue 26 5 '  -0.15%  '

Set-TextValue 27 2 'RenderToken'
Set-TextValue 27 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 27 4 '10.45'
Set-TextValue 27 5 '  -3.48%  '

Set-TextValue 28 4 '3.60'
Set-TextValue 28 5 '  -7.54%  '

Set-TextValue 29 4 '8.98'
Set-TextValue 29 5 '  -5.40%  '

Set-TextValue 30 4 '31.17'
Set-TextValue 30 5 '  -5.37%  '

Set-TextValue 31 4 '6.84'
Set-TextValue 31 5 '  -6.33%  '

Set-TextValue 32 4 '65.98'
Set-TextValue 32 5 '  +2.13%  '

Set-TextValue 33 4 '11.88'
Set-TextValue 33 5 '  -5.08%  '

Set-TextValue 34 2 'Hedera'
Set-TextValue 34 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 34 4 '0.112'
Set-TextValue 34 5 '  -4.58%  '

Set-TextValue 35 2 'InjectiveProtocol'
Set-TextValue 35 3 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 35 4 '42.27'
Set-TextValue 35 5 '  -3.71%  '

Set-TextValue 36 4 '581.03'
Set-TextValue 36 5 '  -3.59%  '

Set-TextValue 37 5 '  -0.08%  '

Set-TextValue 38 5 '  -0.01%  '

Set-TextValue 39 4 '0.376'
Set-TextValue 39 5 '  -6.10%  '

Set-TextValue 40 2 'Kaspa'
Set-TextValue 40 3 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 40 4 '0.133'
Set-TextValue 40 5 '  -2.50%  '

Set-TextValue 41 2 'PEPE'
Set-TextValue 41 3 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 41 4 '0.0₃0727'
Set-TextValue 41 5 '  -18.54%  '

Set-TextValue 42 4 '2.78'
Set-TextValue 42 5 '  -9.04%  '

Set-TextValue 43 4 '0.0411'
Set-TextValue 43 5 '  -6.47%  '

Set-TextValue 44 2 'Stellar'
Set-TextValue 44 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 44 4 '0.132'
Set-TextValue 44 5 '  -2.03%  '

Set-TextValue 45 2 'Fetch.AI'
Set-TextValue 45 3 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 45 4 '2.41'
Set-TextValue 45 5 '  -12.50%  '

Set-TextValue 46 2 'ApeXProtocol'
Set-TextValue 46 3 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 46 4 '3.12'
Set-TextValue 46 5 '  -1.75%  '

Set-TextValue 47 4 '2.693.16'
Set-TextValue 47 5 '  -3.15%  '

Set-TextValue 48 4 '140.83'
Set-TextValue 48 5 '  -0.64%  '

Set-TextValue 49 4 '8.42'
Set-TextValue 49 5 '  -8.93%  '

Set-TextValue 50 4 '2.53'
Set-TextValue 50 5 '  -7.39%  '

Set-TextValue 51 4 '2.58'
Set-TextValue 51 5 '  -7.27%  '
